# Update countries & provincias Spain
#
# Refreshes the COVID "Pais" sheet:
#  - bumps the "Datos actualizados ..." timestamp in A1
#  - refreshes case counters for a handful of existing countries
#  - "Zambia" and "Cabo Verde" climb the (cases-sorted) table: their new,
#    larger totals now slot in above "Mayotte" and "Guinea-Bisau"
#    respectively, pushing the countries between their old and new rank
#    down by one row each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 17 de Julio de 2020 a las 16:12'

# Each entry: row, Pais, Casos totales, Nuevos casos, Casos activos,
#             Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
  @(4,   'Estados Unidos',        3698392, 3367,  1680424, 1876818, 0, 32,  141150),
  @(5,   'Brasil',                2015382, 644,   1366775, 571761,  0, 24,  76846),
  @(6,   'India',                 1017116, 11479, 644172,  347167,  0, 168, 25777),
  @(23,  'Argentina',             114783,  0,     49780,   62870,   0, 21,  2133),
  @(62,  'Serbia',                20109,   392,   14047,   5610,    0, 10,  452),
  @(73,  'Kenia',                 12062,   389,   3983,    7857,    0, 5,   222),
  @(79,  'Noruega',                9018,   3,     8138,    625,     0, 1,   255),
  @(97,  'Republica de Yibuti',    5003,   10,    4809,    138,     0, 0,   56),
  @(109, 'Zambia',                 2810,   190,   1450,    1251,    0, 24,  109),
  @(110, 'Mayotte',                2778,   0,     2581,    160,     0, 0,   37),
  @(111, 'Malaui',                 2712,   0,     1073,    1588,    0, 0,   51),
  @(112, 'Sri Lanka',              2687,   0,     2012,    664,     0, 0,   11),
  @(113, 'Libano',                 2599,   0,     1485,    1074,    0, 0,   40),
  @(114, 'Cuba',                   2444,   4,     2300,    57,      0, 0,   87),
  @(115, 'Mali',                   2440,   0,     1777,    542,     0, 0,   121),
  @(116, 'Congo',                  2358,   0,     589,     1721,    0, 0,   48),
  @(117, 'Sudan del Sur',          2171,   0,     1175,    955,     0, 0,   41),
  @(118, 'Estonia',                2020,   4,     1910,    41,      0, 0,   69),
  @(119, 'Montenegro',             1965,   0,     357,     1582,    0, 0,   26),
  @(120, 'Eslovaquia',             1965,   14,    1523,    414,     0, 0,   28),
  @(121, 'Cabo Verde',             1939,   45,    902,     1018,    0, 0,   19),
  @(122, 'Guinea-Bisau',           1927,   25,    773,     1128,    0, 0,   26),
  @(123, 'Eslovenia',              1916,   19,    1522,    283,     0, 0,   111),
  @(124, 'Islandia',               1916,   2,     1895,    11,      0, 0,   10),
  @(125, 'Lituania',               1908,   6,     1595,    234,     0, 0,   79),
  @(140, 'Liberia',                1085,   15,    496,     520,     0, 1,   69),
  @(147, 'Surinam',                 919,   15,    582,     318,     0, 1,   19),
  @(149, 'Principado de Andorra',   880,   3,     803,     25,      0, 0,   52)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  for ($c = 2; $c -le 8; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c]
  }
}
